$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats enum value
$xlPasteFormats = -4122

# Move the existing data row (row 2) down to row 3: capture its values first
# (must call Value() as a method here, not the bare property, to get the
# real underlying data), then copy its formatting so per-cell styling
# (e.g. the style applied to column A) is preserved on the shifted row.
$oldA = $ws.Range("A2").Value()
$oldB = $ws.Range("B2").Value()
$oldC = $ws.Range("C2").Value()
$oldD = $ws.Range("D2").Value()
$oldE = $ws.Range("E2").Value()
$oldF = $ws.Range("F2").Value()
$oldG = $ws.Range("G2").Value()

$ws.Range("A2:G2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A3").Value = $oldA
$ws.Range("B3").Value = $oldB
$ws.Range("C3").Value = $oldC
$ws.Range("D3").Value = $oldD
$ws.Range("E3").Value = $oldE
$ws.Range("F3").Value = $oldF
$ws.Range("G3").Value = $oldG

# Populate the newly freed row 2 with the new record (formatting stays as
# it already was on row 2, matching the original data-row style pattern).
$ws.Range("A2").Value = "even_MAG-GUT47840.fa"
$ws.Range("B2").Value = 0.359939818383412
$ws.Range("C2").Value = 0.6202911368179811
$ws.Range("D2").Value = 0.01976904479860692
$ws.Range("E2").Value = 0.6202911368179811
$ws.Range("F2").Value = "s__Fenollaria sp900539725"
$ws.Range("G2").Value = "s__Fenollaria sp900539725"
